$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures per latest scrape.
# Force text format so numeric-looking strings (e.g. "157.10", "6.00")
# keep their original formatting instead of becoming Excel numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.362.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.285.22'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '157.10'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +15,606.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '307.10'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '95.44'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +4.56%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.494'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '35.73'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +10.73%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.72'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.639.36'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.46'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.300.46'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.265.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0915'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.91'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '242.57'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.02'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.05'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.82%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.09'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -9.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.93'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.64%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.70%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.21'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.94%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.17'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.010.30'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.38'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +11.00%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.14'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.46%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.37'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.84'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.36%  '
